$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(20, 8).Value = 19304.8
$ws.Cells.Item(20, 9).Value = 750
$ws.Cells.Item(20, 10).Value = 31674.666
$ws.Cells.Item(20, 11).Value = 750
$ws.Cells.Item(20, 12).Value = 31674.666
$ws.Cells.Item(20, 13).Value = -520
$ws.Cells.Item(20, 14).Value = -32134.666
$ws.Cells.Item(35, 8).Value = 19304.8
$ws.Cells.Item(35, 9).Value = 750
$ws.Cells.Item(35, 10).Value = 31674.666
$ws.Cells.Item(35, 11).Value = 750
$ws.Cells.Item(35, 12).Value = 31674.666
$ws.Cells.Item(35, 13).Value = -371
$ws.Cells.Item(35, 14).Value = -32432.666
$ws.Cells.Item(112, 8).Value = 7577100.5
$ws.Cells.Item(112, 9).Value = 2416.6667
$ws.Cells.Item(112, 10).Value = 9260364
$ws.Cells.Item(112, 11).Value = 7250.000100000001
$ws.Cells.Item(112, 12).Value = 27781092
$ws.Cells.Item(112, 13).Value = -6142.000100000001
$ws.Cells.Item(112, 14).Value = -27783308
$ws.Cells.Item(131, 8).Value = 3225.3076
$ws.Cells.Item(131, 10).Value = 2382.25
$ws.Cells.Item(131, 12).Value = 7146.75
$ws.Cells.Item(131, 14).Value = -17226.75
$ws.Cells.Item(132, 8).Value = 8338790
$ws.Cells.Item(132, 9).Value = 9529498
$ws.Cells.Item(132, 10).Value = 3833.3333
$ws.Cells.Item(132, 11).Value = 28588494
$ws.Cells.Item(132, 12).Value = 11499.9999
$ws.Cells.Item(132, 13).Value = -28585964
$ws.Cells.Item(132, 14).Value = -16559.9999
$ws.Cells.Item(138, 8).Value = 3335.4565
$ws.Cells.Item(138, 9).Value = 1544.381
$ws.Cells.Item(138, 10).Value = 4839.96
$ws.Cells.Item(138, 11).Value = 4633.143
$ws.Cells.Item(138, 12).Value = 14519.88
$ws.Cells.Item(138, 13).Value = 506.857
$ws.Cells.Item(138, 14).Value = -24799.88

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 19233960
$ws.Cells.Item(2, 9).Value = 25001246
$ws.Cells.Item(2, 10).Value = 9666.666999999999
$ws.Cells.Item(2, 11).Value = 25001246
$ws.Cells.Item(2, 12).Value = 9666.666999999999
$ws.Cells.Item(2, 13).Value = -25001133
$ws.Cells.Item(2, 14).Value = -9892.666999999999
$ws.Cells.Item(61, 8).Value = 3151.2942
$ws.Cells.Item(61, 9).Value = 1713.1666
$ws.Cells.Item(61, 11).Value = 1713.1666
$ws.Cells.Item(61, 13).Value = -1501.1666
$ws.Cells.Item(116, 8).Value = 19233960
$ws.Cells.Item(116, 9).Value = 25001246
$ws.Cells.Item(116, 10).Value = 9666.666999999999
$ws.Cells.Item(116, 11).Value = 25001246
$ws.Cells.Item(116, 12).Value = 9666.666999999999
$ws.Cells.Item(116, 13).Value = -24998952
$ws.Cells.Item(116, 14).Value = -14254.667
$ws.Cells.Item(122, 8).Value = 2601.9583
$ws.Cells.Item(122, 9).Value = 1731.2667
$ws.Cells.Item(122, 10).Value = 4053.111
$ws.Cells.Item(122, 11).Value = 5193.800099999999
$ws.Cells.Item(122, 12).Value = 12159.333
$ws.Cells.Item(122, 13).Value = -2743.800099999999
$ws.Cells.Item(122, 14).Value = -17059.333
$ws.Cells.Item(132, 8).Value = 22730986
$ws.Cells.Item(132, 9).Value = 32261398
$ws.Cells.Item(132, 10).Value = 4617.385
$ws.Cells.Item(132, 11).Value = 96784194
$ws.Cells.Item(132, 12).Value = 13852.155
$ws.Cells.Item(132, 13).Value = -96781664
$ws.Cells.Item(132, 14).Value = -18912.155
$ws.Cells.Item(136, 8).Value = 3151.2942
$ws.Cells.Item(136, 9).Value = 1713.1666
$ws.Cells.Item(136, 11).Value = 5139.4998
$ws.Cells.Item(136, 13).Value = -2589.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 19233960
$ws.Cells.Item(3, 9).Value = 25001246
$ws.Cells.Item(3, 10).Value = 9666.666999999999
$ws.Cells.Item(3, 11).Value = 25001246
$ws.Cells.Item(3, 12).Value = 9666.666999999999
$ws.Cells.Item(3, 13).Value = -25001132
$ws.Cells.Item(3, 14).Value = -9894.666999999999
$ws.Cells.Item(44, 8).Value = 26350
$ws.Cells.Item(44, 10).Value = 26350
$ws.Cells.Item(44, 12).Value = 26350
$ws.Cells.Item(44, 14).Value = -27344
$ws.Cells.Item(119, 8).Value = 12500
$ws.Cells.Item(119, 10).Value = 12500
$ws.Cells.Item(119, 12).Value = 12500
$ws.Cells.Item(119, 14).Value = -22176

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2634212.8
$ws.Cells.Item(31, 9).Value = 3126740
$ws.Cells.Item(31, 10).Value = 7400
$ws.Cells.Item(31, 11).Value = 3126740
$ws.Cells.Item(31, 12).Value = 7400
$ws.Cells.Item(31, 13).Value = -3126445
$ws.Cells.Item(31, 14).Value = -7990
$ws.Cells.Item(34, 8).Value = 2634212.8
$ws.Cells.Item(34, 9).Value = 3126740
$ws.Cells.Item(34, 10).Value = 7400
$ws.Cells.Item(34, 11).Value = 3126740
$ws.Cells.Item(34, 12).Value = 7400
$ws.Cells.Item(34, 13).Value = -3126538
$ws.Cells.Item(34, 14).Value = -7804
$ws.Cells.Item(99, 8).Value = 2670
$ws.Cells.Item(99, 9).Value = 1172.5
$ws.Cells.Item(99, 11).Value = 1172.5
$ws.Cells.Item(99, 13).Value = 325.5
$ws.Cells.Item(107, 8).Value = 2392.125
$ws.Cells.Item(107, 9).Value = 1546.1538
$ws.Cells.Item(107, 11).Value = 1546.1538
$ws.Cells.Item(107, 13).Value = 373.8462
$ws.Cells.Item(126, 8).Value = 2670
$ws.Cells.Item(126, 9).Value = 1172.5
$ws.Cells.Item(126, 11).Value = 3517.5
$ws.Cells.Item(126, 13).Value = -1047.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 218.73334
$ws.Cells.Item(14, 9).Value = 218.73334
$ws.Cells.Item(14, 11).Value = 656.20002
$ws.Cells.Item(14, 13).Value = -483.20002
$ws.Cells.Item(122, 8).Value = 1959.6
$ws.Cells.Item(122, 10).Value = 2999.3333
$ws.Cells.Item(122, 12).Value = 26993.9997
$ws.Cells.Item(122, 14).Value = -31893.9997
$ws.Cells.Item(131, 8).Value = 897.8472
$ws.Cells.Item(131, 10).Value = 1111.4062
$ws.Cells.Item(131, 12).Value = 3334.2186
$ws.Cells.Item(131, 14).Value = -13414.2186
$ws.Cells.Item(132, 8).Value = 3949.875
$ws.Cells.Item(132, 9).Value = 3899.8
$ws.Cells.Item(132, 10).Value = 4033.3333
$ws.Cells.Item(132, 11).Value = 35098.2
$ws.Cells.Item(132, 12).Value = 36299.9997
$ws.Cells.Item(132, 13).Value = -32568.2
$ws.Cells.Item(132, 14).Value = -41359.9997
$ws.Cells.Item(137, 8).Value = 2338.9412
$ws.Cells.Item(137, 9).Value = 1348.091
$ws.Cells.Item(137, 11).Value = 4044.273
$ws.Cells.Item(137, 13).Value = 1055.727
$ws.Cells.Item(139, 8).Value = 7550.619
$ws.Cells.Item(139, 9).Value = 2420
$ws.Cells.Item(139, 10).Value = 14391.444
$ws.Cells.Item(139, 11).Value = 7260
$ws.Cells.Item(139, 12).Value = 43174.33199999999
$ws.Cells.Item(139, 13).Value = -2120
$ws.Cells.Item(139, 14).Value = -53454.33199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 53.75
$ws.Cells.Item(2, 10).Value = 45.833332
$ws.Cells.Item(2, 12).Value = 45.833332
$ws.Cells.Item(2, 14).Value = -271.833332
$ws.Cells.Item(102, 8).Value = 3415.8333
$ws.Cells.Item(102, 9).Value = 2499
$ws.Cells.Item(102, 11).Value = 2499
$ws.Cells.Item(102, 13).Value = -877

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 3616.889
$ws.Cells.Item(16, 9).Value = 364.57144
$ws.Cells.Item(16, 11).Value = 364.57144
$ws.Cells.Item(16, 13).Value = -194.57144
$ws.Cells.Item(132, 8).Value = 2803.205
$ws.Cells.Item(132, 9).Value = 1832.8182
$ws.Cells.Item(132, 11).Value = 5498.4546
$ws.Cells.Item(132, 13).Value = -2968.4546

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 304426.03
$ws.Cells.Item(132, 9).Value = 481362.28
$ws.Cells.Item(132, 10).Value = 39021.645
$ws.Cells.Item(132, 11).Value = 1444086.84
$ws.Cells.Item(132, 12).Value = 117064.935
$ws.Cells.Item(132, 13).Value = -1441556.84
$ws.Cells.Item(132, 14).Value = -122124.935
$ws.Cells.Item(136, 8).Value = 1797
$ws.Cells.Item(136, 9).Value = 1165.9286
$ws.Cells.Item(136, 10).Value = 2680.5
$ws.Cells.Item(136, 11).Value = 3497.7858
$ws.Cells.Item(136, 12).Value = 8041.5
$ws.Cells.Item(136, 13).Value = -947.7857999999997
$ws.Cells.Item(136, 14).Value = -13141.5
$ws.Cells.Item(138, 8).Value = 29600
$ws.Cells.Item(138, 10).Value = 29600
$ws.Cells.Item(138, 12).Value = 29600
$ws.Cells.Item(138, 14).Value = -39880
$ws.Cells.Item(139, 8).Value = 29490.908
$ws.Cells.Item(139, 10).Value = 29490.908
$ws.Cells.Item(139, 12).Value = 29490.908
$ws.Cells.Item(139, 14).Value = -39770.908
